# Update scripts with new TPM values.
#
# The new TPM run adds an "ECs" sending-cluster row for the Gast -> Cckbr
# edge (targeting FAPs) above the existing FAPs/MuSCs sending-cluster
# rows, and refreshes the numeric specificity/weight columns for all
# three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 20  # columns A..T

# Shift the existing data rows down by one to make room for the new
# "ECs" row: row 3 -> row 4, then row 2 -> row 3. Cell-by-cell value
# copies keep the (unstyled) data-row formatting untouched, unlike a
# Rows.Insert() which would copy the header's formatting down.
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(4, $c).Value = $ws.Cells.Item(3, $c).Value()
}
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(3, $c).Value = $ws.Cells.Item(2, $c).Value()
}

function Set-RowValues($rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $values[$i]
    }
}

# New row 2: Sending cluster "ECs" -> Ligand Gast / Receptor Cckbr -> Target cluster FAPs
Set-RowValues 2 @("ECs", "Gast", "Cckbr", "FAPs", 3, 1, 1.283312333333333, 3.849937, 0.5028409016698769, 0.502840901669877, 3, 1, 0.83582, 2.50746, 1, 1, 1.072618114446666, 9.653563030019999, 0.5028409016698769, 0.502840901669877)

# Row 3 (previously row 2, "FAPs" sending cluster) gets refreshed values
Set-RowValues 3 @("FAPs", "Gast", "Cckbr", "FAPs", 2, 0.6666666666666666, 1.047331333333333, 3.141994, 0.4103763505743974, 0.4103763505743974, 3, 1, 0.83582, 2.50746, 1, 1, 0.8753804750266667, 7.878424275240001, 0.4103763505743974, 0.4103763505743974)

# Row 4 (previously row 3, "MuSCs" sending cluster) gets refreshed values
Set-RowValues 4 @("MuSCs", "Gast", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.2214803333333333, 0.6644409999999999, 0.08678274775572555, 0.08678274775572555, 3, 1, 0.83582, 2.50746, 1, 1, 0.1851176922066667, 1.66605922986, 0.08678274775572555, 0.08678274775572555)
